$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.930.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.326.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.93%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.10'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.71%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.385'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.961'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.326.02'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.56'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.14'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.754.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.949.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.08%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.06'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.323.86'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.29'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.44'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '491.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.440'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -9.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.57'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '90.17'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.92'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.500.18'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.29%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.14'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.90%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.28'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.527'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '561.56'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.43'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.868'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.69'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0415'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.68'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.41'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.11'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.02'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.77%  '
